$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.178.85"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.73%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.574.47"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.29%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "556.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.18%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.57"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.28%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.597"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.19%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.580.83"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.93%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.66"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.89%  "
$ws.Range("E11").Value = "  -1.10%  "
$ws.Range("E12").Value = "  +12.27%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.352"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.50%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.027.38"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.22%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "59.141.51"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.63%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "23.01"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.24%  "
$ws.Range("E17").Value = "  -1.28%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.575.74"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.91%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.55"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.25%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "336.46"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.97%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.33"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.08%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.39"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("E23").Value = "  -0.26%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "62.81"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.20%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.473"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +7.54%  "
$ws.Range("E26").Value = "  +0.41%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.159"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.56%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.41"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0772"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.91%  "
$ws.Range("E30").Value = "  +0.06%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.18"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.29%  "
$ws.Range("E32").Value = "  -3.09%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "157.82"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.44%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.02"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.05%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.03"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.94%  "
$ws.Range("E36").Value = "  -0.11%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.897"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.53%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "37.22"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.53%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.854"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.32%  "
$ws.Range("E40").Value = "  -3.13%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.66"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.07%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "291.03"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.09%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "135.58"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.76%  "
$ws.Range("E44").Value = "  +0.15%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0974"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.12%  "
$ws.Range("E46").Value = "  -2.31%  "
$ws.Range("E47").Value = "  -0.24%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0529"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.99%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0233"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.38%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.60"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.939.81"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.21%  "
